$wb = $excel.ActiveWorkbook

# ALC row 39: Riches' Brew
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 476.86957
$ws.Range("I39").Value = 153.71428
$ws.Range("J39").Value = 979.55554
$ws.Range("K39").Value = 461.14284
$ws.Range("L39").Value = 2938.66662
$ws.Range("M39").Value = -165.14284
$ws.Range("N39").Value = -3530.66662

# ALC row 129: Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 20275.9
$ws.Range("J129").Value = 25296.188
$ws.Range("L129").Value = 75888.564
$ws.Range("N129").Value = -85888.564

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1440.0834
$ws.Range("I137").Value = 1466.7812
$ws.Range("J137").Value = 1226.5
$ws.Range("K137").Value = 4400.3436
$ws.Range("L137").Value = 3679.5
$ws.Range("M137").Value = -1850.3436
$ws.Range("N137").Value = -8779.5

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1371.12
$ws.Range("I138").Value = 678.86206
$ws.Range("J138").Value = 2327.0952
$ws.Range("K138").Value = 2036.58618
$ws.Range("L138").Value = 6981.285600000001
$ws.Range("M138").Value = 3103.41382
$ws.Range("N138").Value = -17261.2856

# ARM row 26: Night Squawker
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1900
$ws.Range("I26").Value = 1900
$ws.Range("K26").Value = 1900
$ws.Range("M26").Value = -1570

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6092.935
$ws.Range("I32").Value = 4754.31
$ws.Range("J32").Value = 21933.334
$ws.Range("K32").Value = 4754.31
$ws.Range("L32").Value = 21933.334
$ws.Range("M32").Value = -4467.31
$ws.Range("N32").Value = -22507.334

# ARM row 68: Let Faith Light the Way
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 28999.25
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 28999.25
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 28999.25
$ws.Range("N68").Value = -30621.25

# ARM row 71: Fifty Shields of Blades (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H71").Value = 28999.25
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 28999.25
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 86997.75
$ws.Range("N71").Value = -95109.75

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 817.8444
$ws.Range("I74").Value = 595.86365
$ws.Range("J74").Value = 1030.174
$ws.Range("K74").Value = 595.86365
$ws.Range("L74").Value = 1030.174
$ws.Range("M74").Value = 278.13635
$ws.Range("N74").Value = -2778.174

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 817.8444
$ws.Range("I77").Value = 595.86365
$ws.Range("J77").Value = 1030.174
$ws.Range("K77").Value = 2979.31825
$ws.Range("L77").Value = 5150.87
$ws.Range("M77").Value = 1388.68175
$ws.Range("N77").Value = -13886.87

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5210174
$ws.Range("I132").Value = 5953770.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 17861311.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -17858781.5
$ws.Range("N132").Value = -20057

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15898485
$ws.Range("I134").Value = 16693260
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 50079780
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -50077245
$ws.Range("N134").Value = -14070

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5750542
$ws.Range("I132").Value = 8773344
$ws.Range("J132").Value = 7218.15
$ws.Range("K132").Value = 26320032
$ws.Range("L132").Value = 21654.45
$ws.Range("M132").Value = -26317502
$ws.Range("N132").Value = -26714.45

# CUL row 2: Pork Is a Salty Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83503.414
$ws.Range("I2").Value = 160
$ws.Range("J2").Value = 100172.1
$ws.Range("K2").Value = 960
$ws.Range("L2").Value = 601032.6000000001
$ws.Range("M2").Value = -847
$ws.Range("N2").Value = -601258.6000000001

# CUL row 4: In Hot Water
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12138.706
$ws.Range("J4").Value = 29248.572
$ws.Range("L4").Value = 87745.716
$ws.Range("N4").Value = -87969.716

# CUL row 5: What a Sap
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 540.381
$ws.Range("I5").Value = 358.16666
$ws.Range("J5").Value = 783.3333
$ws.Range("K5").Value = 1074.49998
$ws.Range("L5").Value = 2349.9999
$ws.Range("M5").Value = -962.4999800000001
$ws.Range("N5").Value = -2573.9999

# CUL row 9: Jack of All Plates
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 250475500
$ws.Range("I9").Value = 1000000000
$ws.Range("J9").Value = 634000
$ws.Range("K9").Value = 3000000000
$ws.Range("L9").Value = 1902000
$ws.Range("M9").Value = -2999999776
$ws.Range("N9").Value = -1902448

# CUL row 10: A Real Fungi
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 923.3333
$ws.Range("I10").Value = 108
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 324
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = -185
$ws.Range("N10").Value = -15278

# CUL row 116: On a Full Stomach
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2499.6667
$ws.Range("I116").Value = 1900
$ws.Range("J116").Value = 2799.5
$ws.Range("K116").Value = 5700
$ws.Range("L116").Value = 8398.5
$ws.Range("M116").Value = -2258
$ws.Range("N116").Value = -15282.5

# CUL row 122: Salt of the North
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 779.2759
$ws.Range("I122").Value = 718.5185
$ws.Range("J122").Value = 1599.5
$ws.Range("K122").Value = 6466.6665
$ws.Range("L122").Value = 14395.5
$ws.Range("M122").Value = -4016.6665
$ws.Range("N122").Value = -19295.5

# CUL row 135: Not-so-secret Ingredient
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 540.381
$ws.Range("I135").Value = 358.16666
$ws.Range("J135").Value = 783.3333
$ws.Range("K135").Value = 3223.49994
$ws.Range("L135").Value = 7049.9997
$ws.Range("M135").Value = -688.4999399999997
$ws.Range("N135").Value = -12119.9997

# GSM row 70: Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 40242.035
$ws.Range("I70").Value = 52003.668
$ws.Range("J70").Value = 4957.143
$ws.Range("K70").Value = 52003.668
$ws.Range("L70").Value = 4957.143
$ws.Range("M70").Value = -51733.668
$ws.Range("N70").Value = -5497.143

# GSM row 73: Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 40242.035
$ws.Range("I73").Value = 52003.668
$ws.Range("J73").Value = 4957.143
$ws.Range("K73").Value = 52003.668
$ws.Range("L73").Value = 4957.143
$ws.Range("M73").Value = -51067.668
$ws.Range("N73").Value = -6829.143

# GSM row 132: On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 20003042
$ws.Range("I132").Value = 25642352
$ws.Range("J132").Value = 9122.182000000001
$ws.Range("K132").Value = 76927056
$ws.Range("L132").Value = 27366.546
$ws.Range("M132").Value = -76924526
$ws.Range("N132").Value = -32426.546

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2382393
$ws.Range("I132").Value = 3334271.2
$ws.Range("J132").Value = 2697.4167
$ws.Range("K132").Value = 10002813.6
$ws.Range("L132").Value = 8092.250100000001
$ws.Range("M132").Value = -10000283.6
$ws.Range("N132").Value = -13152.2501

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10819979
$ws.Range("I132").Value = 6251350
$ws.Range("J132").Value = 15119866
$ws.Range("K132").Value = 18754050
$ws.Range("L132").Value = 45359598
$ws.Range("M132").Value = -18751520
$ws.Range("N132").Value = -45364658

Write-Host "Applied all Fenrir_Profits updates"